#
# Applies the 2019/01/18 "second upload" edit to 3D打印机使用及注意事项.docx
#
# Summary of content changes:
#  1. "...导入SD卡" paragraph gets " 或者 U盘" appended.
#  2. "热床温度和出丝温度（" gets "（平台温度）" inserted after "热床温度".
#  3. "PLA C聚乳酸)" typo fixed to "PLA (聚乳酸)".
#  4. Minor run re-splits (around "stl"/"obj" tokens) that mirror the
#     fine-grained run structure introduced by Word's proofing pass in the
#     original commit (text itself is unchanged there).
#

$d = $word.ActiveDocument

function Split-At($rng) {
    # Forces the engine to break the underlying run at the boundaries of
    # $rng without altering its visible formatting (toggle trick).
    $rng.Bold = 1
    $rng.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "3D模型文件切片，导入SD卡" -> append " 或者 U盘"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("导入SD卡", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $p0 = $rng.Start

    $rng.InsertAfter(" ")
    $p1 = $p0 + 1
    $r1 = $d.Range($p0, $p1)
    $r1.Font.Size = 12

    $rng.InsertAfter("或者")
    $p2 = $p1 + 2
    $r2 = $d.Range($p1, $p2)
    $r2.Font.Size = 12

    $rng.InsertAfter(" U")
    $p3 = $p2 + 2
    $r3 = $d.Range($p2, $p3)
    $r3.Font.Size = 12

    $rng.InsertAfter("盘")
    $p4 = $p3 + 1
    $r4 = $d.Range($p3, $p4)
    $r4.Font.Size = 12
}

# ---------------------------------------------------------------------
# 2) "建立.stl/.obj文" -> split into "建立." / "stl" / "/." / "obj" / "文"
#    (text unchanged, runs re-split to mirror proofing markup)
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("建立.stl/.obj文", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $rng.Start
    Split-At ($d.Range($s + 3, $s + 6))     # "stl"
    Split-At ($d.Range($s + 6, $s + 8))     # "/."
    Split-At ($d.Range($s + 8, $s + 11))    # "obj"
    Split-At ($d.Range($s + 11, $s + 12))   # "文"
}

# ---------------------------------------------------------------------
# 3) "tware里打开建立的.stl/.obj文件" -> same kind of split
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("tware里打开建立的.stl/.obj文件", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $rng.Start
    Split-At ($d.Range($s + 12, $s + 15))   # "stl"
    Split-At ($d.Range($s + 15, $s + 17))   # "/."
    Split-At ($d.Range($s + 17, $s + 20))   # "obj"
    Split-At ($d.Range($s + 20, $s + 22))   # "文件"
}

# ---------------------------------------------------------------------
# 4) "的是热床温度和出丝温度（" -> "的是热床温度" + "（平台温度）" + "和出丝温度（"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("热床温度和出丝温度（", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $insPos = $rng.Start - [int]"和出丝温度（".Length

    $rng2 = $d.Range($insPos, $insPos)
    $rng2.InsertAfter("（平台温度）")

    $newLen = "（平台温度）".Length
    $newRng = $d.Range($insPos, $insPos + $newLen)
    $newRng.Font.Size = 12
    Split-At $newRng

    # also force the boundary split right after the inserted text so the
    # trailing "和出丝温度（" becomes its own run again
    $afterRng = $d.Range($insPos + $newLen, $insPos + $newLen)
    $tailRng = $d.Range($insPos + $newLen, $insPos + $newLen + [int]"和出丝温度（".Length)
    Split-At $tailRng
}

# ---------------------------------------------------------------------
# 5) "PLA C聚乳酸)" -> "PLA (聚乳酸)"
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("PLA C聚乳酸", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $rng.Start
    $cRng = $d.Range($s + 4, $s + 5)
    $cRng.Text = "("
    $cAgain = $d.Range($s + 4, $s + 5)
    Split-At $cAgain
}
